$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers: I0 (col I) and IF (col J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from an existing
# header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-74.
$iVals = @(5,7,8,9,8,8,7,6,7,7,7,6,5,8,7,7,6,7,7,6,7,6,4,7,6,8,5,8,8,4,6,6,6,7,7,8,8,5,8,6,10,6,4,6,6,9,7,8,8,9,7,7,5,6,9,7,6,8,7,6,6,9,6,7,7,5,8,6,5,6,5,5,4)
$jVals = @(5,7,8,9,8,8,7,6,7,7,7,6,5,8,7,7,6,7,7,6,7,6,4,7,6,8,6,8,8,5,6,6,6,7,7,8,9,5,8,6,10,7,5,7,7,9,7,8,9,9,7,7,6,6,9,7,6,8,8,6,6,9,6,7,7,6,8,7,5,6,5,5,4)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
